# Apply the "Updated symbol list" edit (Fri Dec 30 03:39:09 UTC 2022) to Sheet1.
# Every cell in column D holds its numeric reading as literal TEXT (to keep
# trailing zeros / exact digit counts), so every write below is forced to text
# with a leading apostrophe and the cell style is put back to Normal right away
# so no stray number-format style lingers on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, [string]$text)
    $range.Value = "'" + $text
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '245.64'
Set-TextValue $ws.Range("D3") '24.02'
Set-TextValue $ws.Range("D4") '5.252'
Set-TextValue $ws.Range("D5") '0.05784'
Set-TextValue $ws.Range("D6") '6.493'
Set-TextValue $ws.Range("D7") '3.142'
Set-TextValue $ws.Range("D8") '0.8182'
Set-TextValue $ws.Range("D9") '0.8497'
Set-TextValue $ws.Range("D11") '0.06947'
Set-TextValue $ws.Range("D12") '0.03159'
Set-TextValue $ws.Range("D15") '3.745'
Set-TextValue $ws.Range("D16") '0.001514'
Set-TextValue $ws.Range("D17") '0.04717'
Set-TextValue $ws.Range("D18") '0.0005971'
Set-TextValue $ws.Range("D19") '0.006271'
Set-TextValue $ws.Range("D20") '0.001236'
Set-TextValue $ws.Range("D21") '0.004613'
Set-TextValue $ws.Range("E22") '21NitroExNTXWorstin24h'
Set-TextValue $ws.Range("D23") '3.516'
Set-TextValue $ws.Range("D24") '2.133'
Set-TextValue $ws.Range("D27") '0.1326'
Set-TextValue $ws.Range("D28") '0.0002329'
Set-TextValue $ws.Range("D40") '0.03649'
Set-TextValue $ws.Range("B41") 'KickToken'
Set-TextValue $ws.Range("C41") 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
Set-TextValue $ws.Range("D41") '0.006251'
Set-TextValue $ws.Range("E41") '40KickTokenKICK'
Set-TextValue $ws.Range("B42") 'BKEXToken'
Set-TextValue $ws.Range("C42") 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
Set-TextValue $ws.Range("D42") '0.1055'
Set-TextValue $ws.Range("E42") '41BKEXTokenBKK'
Set-TextValue $ws.Range("B43") 'CEJI'
Set-TextValue $ws.Range("C43") 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
Set-TextValue $ws.Range("D43") '0.003400'
Set-TextValue $ws.Range("E43") '42CEJICEJIBestin24h'
Set-TextValue $ws.Range("D44") '0.007460'
Set-TextValue $ws.Range("D45") '0.00005254'
Set-TextValue $ws.Range("D47") '0.3500'
Set-TextValue $ws.Range("D48") '0.002338'
Set-TextValue $ws.Range("D49") '0.00002100'
Set-TextValue $ws.Range("D50") '0.0002000'
